# Applies crypto price/volume updates per commit
# "Updated cryptos list on Thu Sep 14 15:53:35 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.716.23'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '1.637.18'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''213.57'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '''0.492'
$ws.Range('E7').Value = '  +1.05%  '
$ws.Range('E8').Value = '  +0.51%  '
$ws.Range('D9').Value = '''0.0621'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').Value = '''19.08'
$ws.Range('E10').Value = '  +2.98%  '
$ws.Range('E11').Value = '  +2.53%  '
$ws.Range('D12').Value = '1.863.67'
$ws.Range('E12').Value = '  +1.39%  '
$ws.Range('D13').Value = '1.623.07'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').Value = '''4.06'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').Value = '''0.526'
$ws.Range('E15').Value = '  +1.86%  '
$ws.Range('D16').Value = '26.681.94'
$ws.Range('E16').Value = '  +1.26%  '
$ws.Range('D17').Value = '''63.24'
$ws.Range('E17').Value = '  +2.21%  '
$ws.Range('D18').Value = '0.0₃0735'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '''208.86'
$ws.Range('E19').Value = '  +2.64%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '''1.00'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').Value = '''9.39'
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').Value = '''6.10'
$ws.Range('E23').Value = '  +1.04%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').Value = '''146.13'
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  -1.34%  '
$ws.Range('D28').Value = '''15.40'
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('D29').Value = '''6.68'
$ws.Range('E29').Value = '  +1.42%  '
$ws.Range('E30').Value = '  +5.56%  '
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('E32').Value = '  +1.03%  '
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').Value = '''1.51'
$ws.Range('E34').Value = '  +1.08%  '
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').Value = '1.163.47'
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('D38').Value = '''0.813'
$ws.Range('E38').Value = '  +2.08%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '''2.32'
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('D41').Value = '''0.503'
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  +1.00%  '
$ws.Range('D43').Value = '''5.39'
$ws.Range('E43').Value = '  +2.52%  '
$ws.Range('D44').Value = '1.772.91'
$ws.Range('E44').Value = '  +1.22%  '
$ws.Range('D45').Value = '''92.55'
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('D46').Value = '''1.56'
$ws.Range('E46').Value = '  +0.66%  '
$ws.Range('D47').Value = '''54.74'
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('E48').Value = '  +4.57%  '
$ws.Range('E49').Value = '  +0.78%  '
$ws.Range('D50').Value = '''7.66'
$ws.Range('E50').Value = '  +5.17%  '
$ws.Range('D51').Value = '''0.410'
$ws.Range('E51').Value = '  +0.71%  '

# Restore default (unstyled) cell style for numeric-looking text cells
# that Excel would otherwise mark with a quote-prefix style
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
